# Daily attendance processing - 2026-02-01 00:03:14
# Reverse the order of names/emails in the "Recorded By" column (G) for
# every data row where multiple comma-separated entries exist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = $reversed -join ", "
            $cell.Value = $newVal
        }
    }
}
